$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("P2").Value = 1.8
$ws.Range("F3").Value = 5.1
$ws.Range("K3").Value = 13.5
$ws.Range("AF4").Value = 55
$ws.Range("AG4").Value = 30
$ws.Range("AH4").Value = 34
$ws.Range("L5").Value = 1.25
$ws.Range("J6").Value = 3.4
$ws.Range("P6").Value = 1.7
$ws.Range("Q6").Value = 2.12
$ws.Range("I7").Value = 4.2
$ws.Range("J7").Value = 2.66
$ws.Range("L7").Value = 1.45
$ws.Range("M7").Value = 1.11
$ws.Range("P7").Value = 1.57
$ws.Range("Q7").Value = 2.42
$ws.Range("T7").Value = 2
$ws.Range("U7").Value = 1.82
$ws.Range("V7").Value = 1.34
$ws.Range("Q8").Value = 1.99
$ws.Range("H10").Value = 3.95
$ws.Range("F11").Value = 4.1
$ws.Range("F12").Value = 1.71
$ws.Range("I12").Value = 7
$ws.Range("J12").Value = 3
$ws.Range("K12").Value = 4.6
$ws.Range("O12").Value = 1.31
$ws.Range("P12").Value = 1.75
$ws.Range("S12").Value = 3.4
$ws.Range("Y12").Value = 20
$ws.Range("AJ12").Value = 24
$ws.Range("L13").Value = 1.51
$ws.Range("L14").Value = 1.38
$ws.Range("N15").Value = 1.02
$ws.Range("H17").Value = 3.4
$ws.Range("I17").Value = 3.65
$ws.Range("AF18").Value = 16
$ws.Range("J19").Value = 2.96
$ws.Range("L19").Value = 1.47
$ws.Range("P22").Value = 2.16
$ws.Range("H23").Value = 3.35
$ws.Range("L23").Value = 1.38
$ws.Range("Q23").Value = 2.06
$ws.Range("T23").Value = 1.86
$ws.Range("V23").Value = 1.37
$ws.Range("AF23").Value = 17.5
$ws.Range("L24").Value = 1.52
$ws.Range("P24").Value = 1.58
$ws.Range("O25").Value = 1.38
$ws.Range("Q26").Value = 2.32
$ws.Range("G29").Value = 1.75
$ws.Range("Q29").Value = 1.85
$ws.Range("W29").Value = 2.32
$ws.Range("AC29").Value = 9.199999999999999
$ws.Range("AE29").Value = 85
$ws.Range("AG29").Value = 9.800000000000001
$ws.Range("L30").Value = 1.26
$ws.Range("G32").Value = 2.7
$ws.Range("H32").Value = 3.3
$ws.Range("Q32").Value = 2.62
$ws.Range("AN32").Value = 40
$ws.Range("F33").Value = 2.3
$ws.Range("G33").Value = 2.6
$ws.Range("H33").Value = 3
$ws.Range("I33").Value = 3.5
$ws.Range("J33").Value = 3.35
$ws.Range("L33").Value = 1.33
$ws.Range("S33").Value = 3
$ws.Range("T33").Value = 1.72
$ws.Range("U33").Value = 2.14
$ws.Range("W33").Value = 1.63
$ws.Range("Z33").Value = 28
$ws.Range("AA33").Value = 70
$ws.Range("AC33").Value = 9.800000000000001
$ws.Range("AE33").Value = 46
$ws.Range("AF33").Value = 19.5
$ws.Range("AI33").Value = 60
$ws.Range("AO33").Value = 42
$ws.Range("G36").Value = 2.3
$ws.Range("H36").Value = 3.2
$ws.Range("L36").Value = 1.37
$ws.Range("V36").Value = 1.27
$ws.Range("W36").Value = 1.78
$ws.Range("G38").Value = 1.58
$ws.Range("AF39").Value = 16
$ws.Range("F40").Value = 2.5
$ws.Range("V40").Value = 1.44
$ws.Range("F41").Value = 2.5
$ws.Range("V41").Value = 1.4
$ws.Range("G42").Value = 2.86
$ws.Range("H42").Value = 2.7
$ws.Range("I42").Value = 3.05
$ws.Range("J42").Value = 3.35
$ws.Range("W42").Value = 1.54
$ws.Range("F43").Value = 3.4
$ws.Range("G43").Value = 4.1
$ws.Range("I43").Value = 2.3
$ws.Range("L43").Value = 1.2
$ws.Range("N43").Value = 6.6
$ws.Range("O43").Value = 1.14
$ws.Range("R43").Value = 1.78
$ws.Range("T43").Value = 1.46
$ws.Range("U43").Value = 2.74
$ws.Range("V43").Value = 1.77
$ws.Range("W43").Value = 1.34
$ws.Range("F44").Value = 1.47
$ws.Range("L44").Value = 1.18
$ws.Range("F45").Value = 1.79
$ws.Range("H45").Value = 5.1
$ws.Range("N45").Value = 2.22
$ws.Range("O45").Value = 1.65
$ws.Range("S45").Value = 1.02
$ws.Range("U45").Value = 1.53
$ws.Range("T46").Value = 1.96
$ws.Range("K47").Value = 4.2
$ws.Range("P47").Value = 1.76
$ws.Range("AD47").Value = 25
$ws.Range("H48").Value = 1.8
$ws.Range("I48").Value = 1.81
$ws.Range("N48").Value = 3.4
$ws.Range("P48").Value = 1.82
$ws.Range("W48").Value = 1.22
$ws.Range("I50").Value = 6
$ws.Range("P50").Value = 1.58
$ws.Range("V50").Value = 1.2
$ws.Range("G51").Value = 5.1
$ws.Range("H51").Value = 1.73
$ws.Range("I51").Value = 1.74
$ws.Range("R51").Value = 1.48
$ws.Range("U51").Value = 2.28
$ws.Range("V51").Value = 2.34
$ws.Range("W51").Value = 1.24
$ws.Range("X51").Value = 21
$ws.Range("Y51").Value = 9.4
$ws.Range("AD51").Value = 9.4
$ws.Range("AM51").Value = 80
$ws.Range("AO51").Value = 10
$ws.Range("F52").Value = 1.77
$ws.Range("N52").Value = 3.85
$ws.Range("Q52").Value = 2.02
$ws.Range("S52").Value = 3.6
$ws.Range("U52").Value = 1.99
$ws.Range("AG52").Value = 9.6
$ws.Range("AJ52").Value = 19
$ws.Range("F53").Value = 2.86
$ws.Range("G53").Value = 2.96
$ws.Range("H53").Value = 2.88
$ws.Range("I53").Value = 3.05
$ws.Range("N53").Value = 2.56
$ws.Range("P53").Value = 1.53
$ws.Range("AI53").Value = 75
$ws.Range("I54").Value = 2.96
$ws.Range("P54").Value = 1.71
$ws.Range("O55").Value = 1.3
$ws.Range("G56").Value = 2.56
$ws.Range("L57").Value = 1.41
$ws.Range("S57").Value = 4.1
